$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.358.52'
$ws.Range('D3').Value = '1.795.03'
$ws.Range('E3').Value = '  -2.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '307.11'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3598'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.89'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07088'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8847'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07736'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.41'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').Value = '1.788.56'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.282'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.331'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '84.87'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.006'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008518'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.26'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('D22').Value = '26.375.08'
$ws.Range('E22').Value = '  -2.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.981'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.045.14'
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.53'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.967'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.45'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.83'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.023'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '111.96'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.887'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08675'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.058'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.746'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.445'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7251'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.106'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.004'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.067'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('E40').Value = '  -0.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.05086'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.860'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5072'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.65%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.866'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1513'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.006'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.005'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4635'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '101.18'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.23%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.863'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.581'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.15%  '
